$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above the current row 3 (Halloween Event), pushing it
# and the Xmas Party row down to rows 10 and 11. The new rows inherit the
# style of the row above them (row 2: B -> date style, C -> number style).
$ws.Range("A3:A9").EntireRow.Insert()

# The new "menu" rows only use columns A/B (event name + blank date cell);
# drop the C-column cells the insert grew so rows 3-9 keep just A/B.
$ws.Range("C3:C9").Clear()

# Column A: repeat "BBQ Event" down through the new rows. Column B is left
# as the empty, date-styled cell already produced by the row insert.
$ws.Range("A3:A9").Value = "BBQ Event"

# New "Ingredient" / "Member Price" mini-table in columns D/E.
$ws.Range("D1").Value = "Ingredient"
$ws.Range("D2").Value = "Wurst"
$ws.Range("E1").Value = "Member Price"
$ws.Range("D3").Value = "Brochette"
$ws.Range("D4").Value = "Steak"
$ws.Range("D5").Value = "Bearnaise Sauce"
$ws.Range("D6").Value = "Mushroom Sauce"
$ws.Range("D7").Value = "Pepper Sauce"
$ws.Range("D8").Value = "Salmon"
$ws.Range("D9").Value = "Trout"

$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 4
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 5
$ws.Range("E9").Value = 6.05

# Column E uses the same "2 decimal places" number style as column C.
$ws.Range("C1:C2").Copy()
$ws.Range("E1:E9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply the values PasteSpecial's format copy may have clobbered.
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 4
$ws.Range("E4").Value = 6
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 5
$ws.Range("E9").Value = 6.05

# Size the new columns like the other "name"/"value" column pairs (D like
# A, E like C).
$ws.Range("D1").EntireColumn.ColumnWidth = 15.43
$ws.Range("E1").EntireColumn.ColumnWidth = 12.6

# Restore the selection to match the new layout.
$ws.Range("D10").Select()
